# Add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" worksheet (positioned right after "总计" and
#    before "2022-Q2") containing the quarter's fund holdings.
# 2) Insert a new summary row at the top of the "总计" data table for the
#    2022-Q3 totals; existing rows shift down (their own data is untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: build the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$newWs = $wb.Worksheets.Add($beforeSheet)
$newWs.Name = "2022-Q3"

$srcWs = $wb.Worksheets.Item("2022-Q2")

# Copy the header row formatting (bold, centered, bordered = style used by
# every other quarter sheet's header row).
$srcWs.Range("B1:H1").Copy() | Out-Null
$newWs.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# Copy the index-column formatting (column A uses the same bordered style
# down every data row).
$srcWs.Range("A2").Copy() | Out-Null
$newWs.Range("A2:A5").PasteSpecial(-4122) | Out-Null

# Header labels
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Columns B:G hold numeric-looking values that are stored as text (matches
# the source data convention used throughout the workbook).
$newWs.Range("B2:G5").NumberFormat = "@"

# Row 2 - 招商核心竞争力混合A
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "014412"
$newWs.Range("C2").Value = "招商核心竞争力混合A"
$newWs.Range("D2").Value = "22.94"
$newWs.Range("E2").Value = "77.45"
$newWs.Range("F2").Value = "2.87"
$newWs.Range("G2").Value = "0.6584"
$newWs.Range("H2").Value = 7

# Row 3 - 招商核心竞争力混合C
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "014413"
$newWs.Range("C3").Value = "招商核心竞争力混合C"
$newWs.Range("D3").Value = "4.07"
$newWs.Range("E3").Value = "77.45"
$newWs.Range("F3").Value = "2.87"
$newWs.Range("G3").Value = "0.1168"
$newWs.Range("H3").Value = 7

# Row 4 - 招商高端装备混合A
$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "014606"
$newWs.Range("C4").Value = "招商高端装备混合A"
$newWs.Range("D4").Value = "1.56"
$newWs.Range("E4").Value = "93.76"
$newWs.Range("F4").Value = "3.36"
$newWs.Range("G4").Value = "0.0524"
$newWs.Range("H4").Value = 6

# Row 5 - 招商高端装备混合C
$newWs.Range("A5").Value = 3
$newWs.Range("B5").Value = "014607"
$newWs.Range("C5").Value = "招商高端装备混合C"
$newWs.Range("D5").Value = "1.43"
$newWs.Range("E5").Value = "93.76"
$newWs.Range("F5").Value = "3.36"
$newWs.Range("G5").Value = "0.0480"
$newWs.Range("H5").Value = 6

# ---------------------------------------------------------------------
# Part 2: add the 2022-Q3 row to the "总计" (summary) sheet
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Push the existing data rows down by one to make room for the new entry.
$totalWs.Rows(2).Insert()
$totalWs.Range("B2:D2").ClearFormats()

# Re-apply the index-column style (border/bold) that column A carries on
# every data row.
$totalWs.Range("A3").Copy() | Out-Null
$totalWs.Range("A2").PasteSpecial(-4122) | Out-Null

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 4
$totalWs.Range("D2").Value = 0.88

# Renumber the index column (0-based) for the rows that shifted down.
for ($r = 3; $r -le 8; $r++) {
    $totalWs.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q3 data added"
